# Refresh the COVID-19 "paises" data dump: updated timestamp, updated daily
# figures for several countries, and a handful of countries whose relative
# position in the source feed changed (so the row that used to hold one
# country's data now holds another's, and vice versa).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 23:34"

# Update country rows: label text (where the country order changed) and/or
# the statistic columns B..H (where the daily figures changed).

# Row 4
$ws.Cells.Item(4,2).Value = 2202729
$ws.Cells.Item(4,3).Value = 19779
$ws.Cells.Item(4,4).Value = 898077
$ws.Cells.Item(4,5).Value = 1185652
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 717
$ws.Cells.Item(4,8).Value = 119000

# Row 5
$ws.Cells.Item(5,2).Value = 923189
$ws.Cells.Item(5,3).Value = 31633
$ws.Cells.Item(5,4).Value = 464774
$ws.Cells.Item(5,5).Value = 413174
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 1123
$ws.Cells.Item(5,8).Value = 45241

# Row 7
$ws.Cells.Item(7,2).Value = 354161
$ws.Cells.Item(7,3).Value = 11135
$ws.Cells.Item(7,4).Value = 187552
$ws.Cells.Item(7,5).Value = 154688
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 2006
$ws.Cells.Item(7,8).Value = 11921

# Row 13
$ws.Cells.Item(13,2).Value = 188382
$ws.Cells.Item(13,3).Value = 338
$ws.Cells.Item(13,4).Value = 173100
$ws.Cells.Item(13,5).Value = 6372
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 25
$ws.Cells.Item(13,8).Value = 8910

# Row 30
$ws.Cells.Item(30,1).Value = "Ecuador"
$ws.Cells.Item(30,2).Value = 47943
$ws.Cells.Item(30,3).Value = 621
$ws.Cells.Item(30,4).Value = 23684
$ws.Cells.Item(30,5).Value = 20289
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 41
$ws.Cells.Item(30,8).Value = 3970

# Row 31
$ws.Cells.Item(31,1).Value = "Egipto"
$ws.Cells.Item(31,2).Value = 47856
$ws.Cells.Item(31,3).Value = 1567
$ws.Cells.Item(31,4).Value = 12730
$ws.Cells.Item(31,5).Value = 33360
$ws.Cells.Item(31,6).Value = 0
$ws.Cells.Item(31,7).Value = 94
$ws.Cells.Item(31,8).Value = 1766

# Row 49
$ws.Cells.Item(49,1).Value = "Barein"
$ws.Cells.Item(49,2).Value = 19553
$ws.Cells.Item(49,3).Value = 540
$ws.Cells.Item(49,4).Value = 13866
$ws.Cells.Item(49,5).Value = 5640
$ws.Cells.Item(49,6).Value = 0
$ws.Cells.Item(49,7).Value = 1
$ws.Cells.Item(49,8).Value = 47

# Row 50
$ws.Cells.Item(50,1).Value = "Israel"
$ws.Cells.Item(50,2).Value = 19495
$ws.Cells.Item(50,3).Value = 258
$ws.Cells.Item(50,4).Value = 15449
$ws.Cells.Item(50,5).Value = 3744
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 0
$ws.Cells.Item(50,8).Value = 302

# Row 51
$ws.Cells.Item(51,1).Value = "Bolivia"
$ws.Cells.Item(51,2).Value = 19073
$ws.Cells.Item(51,3).Value = 614
$ws.Cells.Item(51,4).Value = 3430
$ws.Cells.Item(51,5).Value = 15011
$ws.Cells.Item(51,6).Value = 0
$ws.Cells.Item(51,7).Value = 21
$ws.Cells.Item(51,8).Value = 632

# Row 90
$ws.Cells.Item(90,2).Value = 3453
$ws.Cells.Item(90,3).Value = 112
$ws.Cells.Item(90,4).Value = 1817
$ws.Cells.Item(90,5).Value = 1455
$ws.Cells.Item(90,6).Value = 0
$ws.Cells.Item(90,7).Value = 5
$ws.Cells.Item(90,8).Value = 181

# Row 126
$ws.Cells.Item(126,1).Value = "Niger"
$ws.Cells.Item(126,2).Value = 1016
$ws.Cells.Item(126,3).Value = 36
$ws.Cells.Item(126,4).Value = 885
$ws.Cells.Item(126,5).Value = 65
$ws.Cells.Item(126,6).Value = 0
$ws.Cells.Item(126,7).Value = 0
$ws.Cells.Item(126,8).Value = 66

# Row 127
$ws.Cells.Item(127,1).Value = "Republica de Chipre"
$ws.Cells.Item(127,2).Value = 985
$ws.Cells.Item(127,3).Value = 0
$ws.Cells.Item(127,4).Value = 807
$ws.Cells.Item(127,5).Value = 160
$ws.Cells.Item(127,6).Value = 0
$ws.Cells.Item(127,7).Value = 0
$ws.Cells.Item(127,8).Value = 18

# Row 128
$ws.Cells.Item(128,1).Value = "Jordania"
$ws.Cells.Item(128,2).Value = 981
$ws.Cells.Item(128,3).Value = 2
$ws.Cells.Item(128,4).Value = 693
$ws.Cells.Item(128,5).Value = 279
$ws.Cells.Item(128,6).Value = 0
$ws.Cells.Item(128,7).Value = 0
$ws.Cells.Item(128,8).Value = 9

# Row 142
$ws.Cells.Item(142,2).Value = 638
$ws.Cells.Item(142,3).Value = 29
$ws.Cells.Item(142,4).Value = 160
$ws.Cells.Item(142,5).Value = 474
$ws.Cells.Item(142,6).Value = 0
$ws.Cells.Item(142,7).Value = 1
$ws.Cells.Item(142,8).Value = 4

# Row 149
$ws.Cells.Item(149,2).Value = 514
$ws.Cells.Item(149,3).Value = 9
$ws.Cells.Item(149,4).Value = 415
$ws.Cells.Item(149,5).Value = 96
$ws.Cells.Item(149,6).Value = 0
$ws.Cells.Item(149,7).Value = 0
$ws.Cells.Item(149,8).Value = 3

# Row 169
$ws.Cells.Item(169,1).Value = "Guyana"
$ws.Cells.Item(169,2).Value = 171
$ws.Cells.Item(169,3).Value = 12
$ws.Cells.Item(169,4).Value = 99
$ws.Cells.Item(169,5).Value = 60
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 0
$ws.Cells.Item(169,8).Value = 12

# Row 170
$ws.Cells.Item(170,1).Value = "Guadalupe"
$ws.Cells.Item(170,2).Value = 171
$ws.Cells.Item(170,3).Value = 0
$ws.Cells.Item(170,4).Value = 157
$ws.Cells.Item(170,5).Value = 0
$ws.Cells.Item(170,6).Value = 0
$ws.Cells.Item(170,7).Value = 0
$ws.Cells.Item(170,8).Value = 14

# Row 206
$ws.Cells.Item(206,1).Value = "Groenlandia"

# Row 207
$ws.Cells.Item(207,1).Value = "Islas Malvinas"

# Row 210
$ws.Cells.Item(210,1).Value = "Seychelles"
$ws.Cells.Item(210,2).Value = 11
$ws.Cells.Item(210,3).Value = 0
$ws.Cells.Item(210,4).Value = 11
$ws.Cells.Item(210,5).Value = 0
$ws.Cells.Item(210,6).Value = 0
$ws.Cells.Item(210,7).Value = 0
$ws.Cells.Item(210,8).Value = 0

# Row 211
$ws.Cells.Item(211,1).Value = "Montserrat"
$ws.Cells.Item(211,2).Value = 11
$ws.Cells.Item(211,3).Value = 0
$ws.Cells.Item(211,4).Value = 10
$ws.Cells.Item(211,5).Value = 0
$ws.Cells.Item(211,6).Value = 0
$ws.Cells.Item(211,7).Value = 0
$ws.Cells.Item(211,8).Value = 1

# Row 213
$ws.Cells.Item(213,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213,2).Value = 8
$ws.Cells.Item(213,3).Value = 0
$ws.Cells.Item(213,4).Value = 8
$ws.Cells.Item(213,5).Value = 0
$ws.Cells.Item(213,6).Value = 0
$ws.Cells.Item(213,7).Value = 0
$ws.Cells.Item(213,8).Value = 0

# Row 214
$ws.Cells.Item(214,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214,2).Value = 8
$ws.Cells.Item(214,3).Value = 0
$ws.Cells.Item(214,4).Value = 7
$ws.Cells.Item(214,5).Value = 0
$ws.Cells.Item(214,6).Value = 0
$ws.Cells.Item(214,7).Value = 0
$ws.Cells.Item(214,8).Value = 1
